$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update "想去人数" (interested count) column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 259
$ws1.Range("F5").Value = 3062
$ws1.Range("F6").Value = 2048
$ws1.Range("F8").Value = 143
$ws1.Range("F9").Value = 1142
$ws1.Range("F10").Value = 209
$ws1.Range("F11").Value = 849

# Sheet "全部类型" (All types) - update "想去人数" (interested count) column F
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 259
$ws4.Range("F5").Value = 3062
$ws4.Range("F6").Value = 2048
$ws4.Range("F9").Value = 143
$ws4.Range("F10").Value = 1142
$ws4.Range("F11").Value = 209
$ws4.Range("F12").Value = 849
